$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.115.35"
$ws.Cells.Item(2, 5).Value = "  +3.40%  "
$ws.Cells.Item(3, 4).Value = "2.986.12"
$ws.Cells.Item(3, 5).Value = "  +2.40%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.05%  "
$ws.Cells.Item(5, 4).Value = "'597.33"
$ws.Cells.Item(5, 5).Value = "  +1.17%  "
$ws.Cells.Item(6, 4).Value = "'146.62"
$ws.Cells.Item(6, 5).Value = "  +1.07%  "
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 4).Value = "2.980.09"
$ws.Cells.Item(8, 5).Value = "  +2.26%  "
$ws.Cells.Item(9, 4).Value = "'0.505"
$ws.Cells.Item(9, 5).Value = "  +0.12%  "
$ws.Cells.Item(10, 4).Value = "'7.42"
$ws.Cells.Item(10, 5).Value = "  +6.93%  "
$ws.Cells.Item(11, 4).Value = "'0.145"
$ws.Cells.Item(11, 5).Value = "  +2.81%  "
$ws.Cells.Item(12, 4).Value = "'0.452"
$ws.Cells.Item(12, 5).Value = "  +3.10%  "
$ws.Cells.Item(13, 4).Value = "'0.0000237"
$ws.Cells.Item(13, 5).Value = "  +5.68%  "
$ws.Cells.Item(14, 4).Value = "'33.72"
$ws.Cells.Item(14, 5).Value = "  +0.89%  "
$ws.Cells.Item(15, 5).Value = "  +0.48%  "
$ws.Cells.Item(16, 4).Value = "3.480.57"
$ws.Cells.Item(16, 5).Value = "  +2.42%  "
$ws.Cells.Item(17, 4).Value = "62.867.11"
$ws.Cells.Item(17, 5).Value = "  +3.17%  "
$ws.Cells.Item(18, 4).Value = "'6.77"
$ws.Cells.Item(18, 5).Value = "  +1.39%  "
$ws.Cells.Item(19, 4).Value = "2.963.82"
$ws.Cells.Item(19, 5).Value = "  +1.61%  "
$ws.Cells.Item(20, 4).Value = "'444.92"
$ws.Cells.Item(20, 5).Value = "  +2.66%  "
$ws.Cells.Item(21, 4).Value = "'13.61"
$ws.Cells.Item(21, 5).Value = "  +1.87%  "
$ws.Cells.Item(22, 4).Value = "'0.679"
$ws.Cells.Item(22, 5).Value = "  +0.57%  "
$ws.Cells.Item(23, 4).Value = "'7.18"
$ws.Cells.Item(23, 5).Value = "  +1.09%  "
$ws.Cells.Item(24, 4).Value = "'82.38"
$ws.Cells.Item(24, 5).Value = "  +1.21%  "
$ws.Cells.Item(25, 4).Value = "'10.93"
$ws.Cells.Item(25, 5).Value = "  +0.73%  "
$ws.Cells.Item(26, 4).Value = "'12.17"
$ws.Cells.Item(26, 5).Value = "  +3.30%  "
$ws.Cells.Item(27, 4).Value = "'2.16"
$ws.Cells.Item(27, 5).Value = "  -2.14%  "
$ws.Cells.Item(28, 5).Value = "  +0.06%  "
$ws.Cells.Item(29, 4).Value = "'2.63"
$ws.Cells.Item(29, 5).Value = "  +1.51%  "
$ws.Cells.Item(30, 4).Value = "'7.18"
$ws.Cells.Item(30, 5).Value = "  +3.03%  "
$ws.Cells.Item(31, 4).Value = "'2.14"
$ws.Cells.Item(31, 5).Value = "  -5.91%  "
$ws.Cells.Item(32, 4).Value = "'26.66"
$ws.Cells.Item(32, 5).Value = "  +0.72%  "
$ws.Cells.Item(33, 4).Value = "'0.109"
$ws.Cells.Item(33, 5).Value = "  +0.36%  "
$ws.Cells.Item(34, 4).Value = "'0.999"
$ws.Cells.Item(34, 5).Value = "  -0.07%  "
$ws.Cells.Item(35, 4).Value = "0.0₃0886"
$ws.Cells.Item(35, 5).Value = "  +2.01%  "
$ws.Cells.Item(36, 4).Value = "'0.996"
$ws.Cells.Item(36, 5).Value = "  -1.58%  "
$ws.Cells.Item(37, 4).Value = "'5.68"
$ws.Cells.Item(37, 5).Value = "  +1.24%  "
$ws.Cells.Item(38, 4).Value = "'2.06"
$ws.Cells.Item(38, 5).Value = "  +4.65%  "
$ws.Cells.Item(39, 2).Value = "dogwifhat"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(39, 4).Value = "'3.00"
$ws.Cells.Item(39, 5).Value = "  -0.20%  "
$ws.Cells.Item(40, 2).Value = "OKB"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(40, 4).Value = "'49.81"
$ws.Cells.Item(40, 5).Value = "  +0.47%  "
$ws.Cells.Item(41, 4).Value = "'8.66"
$ws.Cells.Item(41, 5).Value = "  +1.21%  "
$ws.Cells.Item(42, 4).Value = "'0.119"
$ws.Cells.Item(42, 5).Value = "  -2.16%  "
$ws.Cells.Item(43, 4).Value = "'0.288"
$ws.Cells.Item(43, 5).Value = "  +0.28%  "
$ws.Cells.Item(44, 4).Value = "'39.09"
$ws.Cells.Item(44, 5).Value = "  -6.52%  "
$ws.Cells.Item(45, 4).Value = "'373.85"
$ws.Cells.Item(45, 5).Value = "  -0.51%  "
$ws.Cells.Item(46, 2).Value = "VeChain"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(46, 4).Value = "'0.0345"
$ws.Cells.Item(46, 5).Value = "  -0.23%  "
$ws.Cells.Item(47, 2).Value = "Maker"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(47, 4).Value = "2.708.72"
$ws.Cells.Item(47, 5).Value = "  +0.66%  "
$ws.Cells.Item(48, 4).Value = "'134.27"
$ws.Cells.Item(48, 5).Value = "  +0.97%  "
$ws.Cells.Item(49, 5).Value = "  +0.03%  "
$ws.Cells.Item(50, 4).Value = "'23.52"
$ws.Cells.Item(50, 5).Value = "  -1.20%  "
$ws.Cells.Item(51, 4).Value = "'0.106"
$ws.Cells.Item(51, 5).Value = "  +0.06%  "

# Remove quote-prefix formatting introduced by forcing numeric-looking text
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(51, 4).ClearFormats()
